$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Kode_PKS (A2) with new code. Use a leading apostrophe to keep the
# "number stored as text" quote-prefix formatting that the cell already had.
$ws.Range("A2").Value = "'01732290"

# Update the active selection on the sheet
$ws.Range("D4").Select()

# Update the workbook window size/position
$win = $wb.Windows.Item(1)
$win.Left = 1152
$win.Top = 1152
$win.Width = 15972
$win.Height = 10764
